$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") changed from 2023-09-23 (45192) to 2023-10-03 (45202)
#    for every existing data row (rows 2-205).
for ($r = 2; $r -le 205; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2) Row 205 gains an explicit default row height (ht="15" customHeight="1"),
#    matching the other data rows.
$ws.Rows.Item(205).RowHeight = 15

# 3) A brand-new row 206 is appended with a new case.
$ws.Cells.Item(206, 1).Value = "A 45588-2023"

$ws.Cells.Item(206, 2).Value = 45194
$ws.Cells.Item(206, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(206, 3).Value = 45202
$ws.Cells.Item(206, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(206, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(206, 5).Value = "ÖRKELLJUNGA"
$ws.Cells.Item(206, 6).Value = "Kyrkan"

$ws.Cells.Item(206, 7).Value = 1.5
$ws.Cells.Item(206, 8).Value = 0
$ws.Cells.Item(206, 9).Value = 0
$ws.Cells.Item(206, 10).Value = 0
$ws.Cells.Item(206, 11).Value = 0
$ws.Cells.Item(206, 12).Value = 0
$ws.Cells.Item(206, 13).Value = 0
$ws.Cells.Item(206, 14).Value = 0
$ws.Cells.Item(206, 15).Value = 0
$ws.Cells.Item(206, 16).Value = 0
$ws.Cells.Item(206, 17).Value = 0

$ws.Cells.Item(206, 18).Value = ""
$ws.Cells.Item(206, 18).WrapText = $true
